# Apply updated cryptos list values (price + 1h volume change) per commit.
# Values that look like plain numbers are prefixed with a leading apostrophe
# so Excel stores them as text (matching the sheet's existing text-formatted
# Price/Volume columns) instead of silently converting them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.947.10'
$ws.Range('E2').Value = '  +0.51%  '
$ws.Range('D3').Value = '1.650.17'
$ws.Range('E3').Value = '  +1.04%  '
$ws.Range('D4').Value = '''1.009'
$ws.Range('E4').Value = '  +0.81%  '
$ws.Range('D5').Value = '''216.42'
$ws.Range('E5').Value = '  +0.78%  '
$ws.Range('D6').Value = '''0.5111'
$ws.Range('E6').Value = '  +1.74%  '
$ws.Range('D7').Value = '''1.009'
$ws.Range('E7').Value = '  +0.74%  '
$ws.Range('D8').Value = '''0.2585'
$ws.Range('E8').Value = '  +0.84%  '
$ws.Range('D9').Value = '''0.06439'
$ws.Range('E9').Value = '  +0.91%  '
$ws.Range('E10').Value = '  +0.38%  '
$ws.Range('D11').Value = '''0.07802'
$ws.Range('E11').Value = '  +1.54%  '
$ws.Range('D12').Value = '''4.331'
$ws.Range('E12').Value = '  +2.11%  '
$ws.Range('D13').Value = '1.651.48'
$ws.Range('E13').Value = '  +1.08%  '
$ws.Range('D14').Value = '''0.5480'
$ws.Range('E14').Value = '  +1.09%  '
$ws.Range('D15').Value = '0.0₅7906'
$ws.Range('E15').Value = '  +0.04%  '
$ws.Range('D16').Value = '''64.88'
$ws.Range('E16').Value = '  +2.19%  '
$ws.Range('D17').Value = '26.051.37'
$ws.Range('E17').Value = '  +0.87%  '
$ws.Range('D18').Value = '''1.009'
$ws.Range('E18').Value = '  +0.74%  '
$ws.Range('D19').Value = '''198.70'
$ws.Range('E19').Value = '  -1.63%  '
$ws.Range('D20').Value = '''4.487'
$ws.Range('E20').Value = '  +3.69%  '
$ws.Range('D21').Value = '''10.05'
$ws.Range('E21').Value = '  +1.31%  '
$ws.Range('D22').Value = '''6.089'
$ws.Range('E22').Value = '  +2.24%  '
$ws.Range('D23').Value = '''1.011'
$ws.Range('E23').Value = '  +0.82%  '
$ws.Range('D24').Value = '''1.865'
$ws.Range('E24').Value = '  -3.59%  '
$ws.Range('D25').Value = '''140.60'
$ws.Range('E25').Value = '  -0.57%  '
$ws.Range('D26').Value = '''0.1154'
$ws.Range('E26').Value = '  +1.23%  '
$ws.Range('D27').Value = '''6.920'
$ws.Range('E27').Value = '  +3.28%  '
$ws.Range('D28').Value = '''15.77'
$ws.Range('E28').Value = '  +0.66%  '
$ws.Range('D29').Value = '''1.246'
$ws.Range('E29').Value = '  +0.66%  '
$ws.Range('D30').Value = '''0.05030'
$ws.Range('E30').Value = '  +0.72%  '
$ws.Range('E31').Value = '  +1.01%  '
$ws.Range('D32').Value = '''3.212'
$ws.Range('E32').Value = '  +1.09%  '
$ws.Range('D33').Value = '''1.550'
$ws.Range('E33').Value = '  +0.79%  '
$ws.Range('D34').Value = '''2.373'
$ws.Range('E34').Value = '  +0.38%  '
$ws.Range('D35').Value = '''0.8965'
$ws.Range('E35').Value = '  +0.58%  '
$ws.Range('D36').Value = '''2.604'
$ws.Range('E36').Value = '  -0.32%  '
$ws.Range('D37').Value = '1.138.19'
$ws.Range('E37').Value = '  -2.78%  '
$ws.Range('D38').Value = '''0.5548'
$ws.Range('D39').Value = '''0.01569'
$ws.Range('E39').Value = '  +0.64%  '
$ws.Range('D40').Value = '''1.011'
$ws.Range('E40').Value = '  +1.00%  '
$ws.Range('D41').Value = '''5.678'
$ws.Range('E41').Value = '  -0.03%  '
$ws.Range('D42').Value = '''0.8203'
$ws.Range('E42').Value = '  +1.76%  '
$ws.Range('B43').Value = 'BabyDogeCoin'
$ws.Range('C43').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D43').Value = '0.0₈126'
$ws.Range('E43').Value = '  +9.60%  '
$ws.Range('B44').Value = 'Quant'
$ws.Range('C44').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D44').Value = '''100.17'
$ws.Range('E44').Value = '  +0.83%  '
$ws.Range('D45').Value = '1.785.69'
$ws.Range('E45').Value = '  +0.89%  '
$ws.Range('D46').Value = '''0.4546'
$ws.Range('E46').Value = '  +0.76%  '
$ws.Range('D47').Value = '''55.40'
$ws.Range('E47').Value = '  +1.38%  '
$ws.Range('D49').Value = '''0.05099'
$ws.Range('E49').Value = '  +0.39%  '
$ws.Range('D50').Value = '''1.010'
$ws.Range('E50').Value = '  +0.78%  '
$ws.Range('D51').Value = '''0.09561'
$ws.Range('E51').Value = '  +3.29%  '
